$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Value = 2.62
$ws.Range("K2").Value = 1.87
$ws.Range("L3").Value = 7.5
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("Z3").Value = 9.5
$ws.Range("AE3").Value = 26
$ws.Range("AF3").Value = 101
$ws.Range("AZ3").Value = 201
$ws.Range("BA3").Value = 251
$ws.Range("C5").Value = "18:30"
$ws.Range("G5").Value = 1.3
$ws.Range("H5").Value = 5.25
$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 1.8
$ws.Range("K5").Value = 2.5
$ws.Range("L5").Value = 9
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 2.08
$ws.Range("U5").Value = 2.2
$ws.Range("V5").Value = 1.62
$ws.Range("N8").Value = 8.5
$ws.Range("Y8").Value = 13
$ws.Range("AN8").Value = 5.5
$ws.Range("AS8").Value = 251
$ws.Range("H11").Value = 3.75
$ws.Range("K11").Value = 2.1
$ws.Range("N11").Value = 8.5
$ws.Range("W11").Value = 6
$ws.Range("AA11").Value = 15
$ws.Range("AC11").Value = 8.5
$ws.Range("AG11").Value = 11
$ws.Range("AQ11").Value = 29
$ws.Range("BB11").Value = 351
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
